$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("grandes regiões" label with no data) is removed; the rows below
# (norte, nordeste, sudeste, sul with their values) shift up by one.
$ws.Rows.Item(6).Delete()
